# Generate Report for Handoff
# Replaces the two tracked source files (and their associated handoff/handback
# metadata) with a new pair of files that are "Ready for handoff":
#   4716211e-49fa-4399-b42a-270874a094ef.md -> ef67d290-90fe-4bad-9428-55186c0e639b.md
#   7cc93c55-f71b-41a4-8a8e-7cafe7a27353.md -> ffff18aa7941-0904-4bb0-92b6-ab27b4ff6532.md
# Status changes from "Handed back: in sync with en-US" to "Ready for handoff",
# the handback related columns (Latest Target File / Latest Handback File) are
# dropped, and the handoff file/datetime columns are refreshed.

$wb = $excel.ActiveWorkbook

$oldFile1 = "4716211e-49fa-4399-b42a-270874a094ef.md"
$oldFile2 = "7cc93c55-f71b-41a4-8a8e-7cafe7a27353.md"

$newFile1 = "ef67d290-90fe-4bad-9428-55186c0e639b.md"
$newFile2 = "ffff18aa7941-0904-4bb0-92b6-ab27b4ff6532.md"

$newStatus = "Ready for handoff"

$newXlfZh = "ef67d290-90fe-4bad-9428-55186c0e639b.4999ebf1b2ed29a4b1220a5efa5288b599a1e03f.zh-cn.xlf"
$newXlfDe = "ef67d290-90fe-4bad-9428-55186c0e639b.4999ebf1b2ed29a4b1220a5efa5288b599a1e03f.de-de.xlf"

$newHandoffDtZh = "2016-03-10 17:00:11"
$newHandoffDtDe = "2016-03-10 17:00:21"
$noTarget = "0001-01-01 00:00:00"

function Set-HyperlinkInPlace($ws, $addr, $displayText, $targetUrl) {
    $all = @($ws.Hyperlinks)
    foreach ($hl in $all) {
        if ($hl.Range.Address() -eq $addr) {
            $hl.TextToDisplay = $displayText
            $hl.Address = $targetUrl
            return
        }
    }
}

function Remove-HyperlinkAt($ws, $addr) {
    $all = @($ws.Hyperlinks)
    foreach ($hl in $all) {
        if ($hl.Range.Address() -eq $addr) {
            $hl.Delete()
            return
        }
    }
}

# ---------------------------------------------------------------------------
# Sheet "Overview"
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("A2").Value = $newFile1
$wsOverview.Range("B2").Value = $newStatus
$wsOverview.Range("C2").Value = $newStatus

$wsOverview.Range("A3").Value = $newFile2
$wsOverview.Range("B3").Value = $newStatus
$wsOverview.Range("C3").Value = $newStatus

Set-HyperlinkInPlace $wsOverview "`$A`$2" $newFile1 "https://github.com/OpenLocalizationTest/oltest/blob/7d224905e22ba47e275323ab53d06186e4677670/e2e/$newFile1"
Set-HyperlinkInPlace $wsOverview "`$A`$3" $newFile2 "https://github.com/OpenLocalizationTest/oltest/blob/7d224905e22ba47e275323ab53d06186e4677670/e2e/$newFile2"

# ---------------------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Range("A2").Value = $newFile1
$wsZh.Range("B2").Value = $newStatus
$wsZh.Range("C2").Value = $newXlfZh
$wsZh.Range("D2").Value = $newHandoffDtZh
$wsZh.Range("G2").Value = $noTarget

$wsZh.Range("A3").Value = $newFile2
$wsZh.Range("B3").Value = $newStatus
$wsZh.Range("C3").Value = $newXlfZh
$wsZh.Range("D3").Value = $newHandoffDtZh
$wsZh.Range("G3").Value = $noTarget

# Drop the "Latest Target File" / "Latest Handback File" columns for rows 2-3
$wsZh.Range("E2:F3").Clear()

Remove-HyperlinkAt $wsZh "`$E`$2"
Remove-HyperlinkAt $wsZh "`$F`$2"
Remove-HyperlinkAt $wsZh "`$E`$3"
Remove-HyperlinkAt $wsZh "`$F`$3"

Set-HyperlinkInPlace $wsZh "`$A`$2" $newFile1 "https://github.com/OpenLocalizationTest/oltest/blob/7d224905e22ba47e275323ab53d06186e4677670/e2e/$newFile1"
Set-HyperlinkInPlace $wsZh "`$C`$2" $newXlfZh "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/b2fb6902910a7fbf4713cccb201eb3273d12f454/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/$newXlfZh"
Set-HyperlinkInPlace $wsZh "`$A`$3" $newFile2 "https://github.com/OpenLocalizationTest/oltest/blob/7d224905e22ba47e275323ab53d06186e4677670/e2e/$newFile2"
Set-HyperlinkInPlace $wsZh "`$C`$3" $newXlfZh "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/b2fb6902910a7fbf4713cccb201eb3273d12f454/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/$newXlfZh"

# ---------------------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Range("A2").Value = $newFile1
$wsDe.Range("B2").Value = $newStatus
$wsDe.Range("C2").Value = $newXlfDe
$wsDe.Range("D2").Value = $newHandoffDtDe
$wsDe.Range("G2").Value = $noTarget

$wsDe.Range("A3").Value = $newFile2
$wsDe.Range("B3").Value = $newStatus
$wsDe.Range("C3").Value = $newXlfDe
$wsDe.Range("D3").Value = $newHandoffDtDe
$wsDe.Range("G3").Value = $noTarget

# Drop the "Latest Target File" / "Latest Handback File" columns for rows 2-3
$wsDe.Range("E2:F3").Clear()

Remove-HyperlinkAt $wsDe "`$E`$2"
Remove-HyperlinkAt $wsDe "`$F`$2"
Remove-HyperlinkAt $wsDe "`$E`$3"
Remove-HyperlinkAt $wsDe "`$F`$3"

Set-HyperlinkInPlace $wsDe "`$A`$2" $newFile1 "https://github.com/OpenLocalizationTest/oltest/blob/7d224905e22ba47e275323ab53d06186e4677670/e2e/$newFile1"
Set-HyperlinkInPlace $wsDe "`$C`$2" $newXlfDe "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/a1201cd6e701b8ccfb4922dbab86dded67ebab93/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/$newXlfDe"
Set-HyperlinkInPlace $wsDe "`$A`$3" $newFile2 "https://github.com/OpenLocalizationTest/oltest/blob/7d224905e22ba47e275323ab53d06186e4677670/e2e/$newFile2"
Set-HyperlinkInPlace $wsDe "`$C`$3" $newXlfDe "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/a1201cd6e701b8ccfb4922dbab86dded67ebab93/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/$newXlfDe"
